$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("E2").Value = 3
$ws.Range("F2").Value = 1
$ws.Range("G2").Value = 3.537653
$ws.Range("H2").Value = 10.612959
$ws.Range("I2").Value = 0.277393541355334
$ws.Range("J2").Value = 0.277393541355334
$ws.Range("M2").Value = 0.09834766666666667
$ws.Range("N2").Value = 0.295043
$ws.Range("O2").Value = 0.2818566198948398
$ws.Range("P2").Value = 0.2818566198948398
$ws.Range("Q2").Value = 0.3479199180263334
$ws.Range("R2").Value = 3.131279262237
$ws.Range("S2").Value = 0.07818520594707389
$ws.Range("T2").Value = 0.0781852059470739
$ws.Range("E3").Value = 3
$ws.Range("F3").Value = 1
$ws.Range("G3").Value = 3.537653
$ws.Range("H3").Value = 10.612959
$ws.Range("I3").Value = 0.277393541355334
$ws.Range("J3").Value = 0.277393541355334
$ws.Range("K3").Value = 2
$ws.Range("L3").Value = 0.6666666666666666
$ws.Range("M3").Value = 0.2505803333333333
$ws.Range("N3").Value = 0.751741
$ws.Range("O3").Value = 0.7181433801051602
$ws.Range("P3").Value = 0.7181433801051602
$ws.Range("Q3").Value = 0.8864662679576668
$ws.Range("R3").Value = 7.978196411619
$ws.Range("S3").Value = 0.1992083354082601
$ws.Range("T3").Value = 0.1992083354082601
$ws.Range("I4").Value = 0.422881095777985
$ws.Range("J4").Value = 0.422881095777985
$ws.Range("M4").Value = 0.09834766666666667
$ws.Range("N4").Value = 0.295043
$ws.Range("O4").Value = 0.2818566198948398
$ws.Range("P4").Value = 0.2818566198948398
$ws.Range("Q4").Value = 0.5303971947547778
$ws.Range("R4").Value = 4.773574752793
$ws.Range("S4").Value = 0.1191918362734089
$ws.Range("T4").Value = 0.1191918362734089
$ws.Range("I5").Value = 0.422881095777985
$ws.Range("J5").Value = 0.422881095777985
$ws.Range("K5").Value = 2
$ws.Range("L5").Value = 0.6666666666666666
$ws.Range("M5").Value = 0.2505803333333333
$ws.Range("N5").Value = 0.751741
$ws.Range("O5").Value = 0.7181433801051602
$ws.Range("P5").Value = 0.7181433801051602
$ws.Range("Q5").Value = 1.351400702887889
$ws.Range("R5").Value = 12.162606325991
$ws.Range("S5").Value = 0.3036892595045761
$ws.Range("T5").Value = 0.3036892595045762
$ws.Range("G6").Value = 3.797319666666667
$ws.Range("H6").Value = 11.391959
$ws.Range("I6").Value = 0.2977544575442879
$ws.Range("J6").Value = 0.2977544575442879
$ws.Range("M6").Value = 0.09834766666666667
$ws.Range("N6").Value = 0.295043
$ws.Range("O6").Value = 0.2818566198948398
$ws.Range("P6").Value = 0.2818566198948398
$ws.Range("Q6").Value = 0.3734575288041111
$ws.Range("R6").Value = 3.361117759237
$ws.Range("S6").Value = 0.08392406496205457
$ws.Range("T6").Value = 0.08392406496205458
$ws.Range("G7").Value = 3.797319666666667
$ws.Range("H7").Value = 11.391959
$ws.Range("I7").Value = 0.2977544575442879
$ws.Range("J7").Value = 0.2977544575442879
$ws.Range("K7").Value = 2
$ws.Range("L7").Value = 0.6666666666666666
$ws.Range("M7").Value = 0.2505803333333333
$ws.Range("N7").Value = 0.751741
$ws.Range("O7").Value = 0.7181433801051602
$ws.Range("P7").Value = 0.7181433801051602
$ws.Range("Q7").Value = 0.9515336278465556
$ws.Range("R7").Value = 8.563802650618999
$ws.Range("S7").Value = 0.2138303925822334
$ws.Range("T7").Value = 0.2138303925822334
$ws.Range("G8").Value = 0.02513533333333333
$ws.Range("H8").Value = 0.075406
$ws.Range("I8").Value = 0.001970905322393153
$ws.Range("J8").Value = 0.001970905322393153
$ws.Range("M8").Value = 0.09834766666666667
$ws.Range("N8").Value = 0.295043
$ws.Range("O8").Value = 0.2818566198948398
$ws.Range("P8").Value = 0.2818566198948398
$ws.Range("Q8").Value = 0.002472001384222222
$ws.Range("R8").Value = 0.022248012458
$ws.Range("S8").Value = 0.0005555127123024836
$ws.Range("T8").Value = 0.0005555127123024837
$ws.Range("G9").Value = 0.02513533333333333
$ws.Range("H9").Value = 0.075406
$ws.Range("I9").Value = 0.001970905322393153
$ws.Range("J9").Value = 0.001970905322393153
$ws.Range("K9").Value = 2
$ws.Range("L9").Value = 0.6666666666666666
$ws.Range("M9").Value = 0.2505803333333333
$ws.Range("N9").Value = 0.751741
$ws.Range("O9").Value = 0.7181433801051602
$ws.Range("P9").Value = 0.7181433801051602
$ws.Range("Q9").Value = 0.006298420205111111
$ws.Range("R9").Value = 0.056685781846
$ws.Range("S9").Value = 0.001415392610090669
$ws.Range("T9").Value = 0.001415392610090669
